$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new test row (row 85) matching the pattern of the other rows:
# Column A = Test name, Column B = Description, Column C = macro name
# (values set in this order so new shared strings are created in the
# same order as the target workbook)
$ws.Cells.Item(85, 3).Value = "Drawdown_Peak_test1"
$ws.Cells.Item(85, 1).Value = "Drawdown Peak1"
$ws.Cells.Item(85, 2).Value = "Test drawdown_peak"

# Update the selection to match the post-edit state
$ws.Range("C86").Select()
